$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are always stored as text in this sheet (e.g. thousand
# separators rendered as extra dots). A leading apostrophe forces Excel to treat
# the assigned value as text instead of auto-converting it to a number.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'44.833.54"
$ws.Range("E2").Value = "  +1.67%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.250.26"
$ws.Range("E3").Value = "  +0.44%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.27%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'306.86"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'96.00"
$ws.Range("E6").Value = "  +0.01%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - USDC
$ws.Range("D8").Value = "'1.01"
$ws.Range("E8").Value = "  +0.15%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.76%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'35.41"
$ws.Range("E10").Value = "  +1.52%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.68%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +0.10%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.03%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'2.592.58"
$ws.Range("E14").Value = "  +0.24%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "'2.308.72"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "'0.844"
$ws.Range("E16").Value = "  +2.02%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "'13.62"
$ws.Range("E17").Value = "  +0.37%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'44.649.70"
$ws.Range("E18").Value = "  +1.43%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "'0.0₃0951"
$ws.Range("E19").Value = "  -1.57%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.06"
$ws.Range("E20").Value = "  -1.25%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'6.29"
$ws.Range("E21").Value = "  -0.78%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'65.55"
$ws.Range("E22").Value = "  +0.38%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'241.02"
$ws.Range("E23").Value = "  +1.96%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  +0.57%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -0.05%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.02%  "

# Row 27 - Toncoin
$ws.Range("D27").Value = "'2.29"
$ws.Range("E27").Value = "  +3.40%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "'9.90"
$ws.Range("E28").Value = "  +0.03%  "

# Row 29 - InjectiveProtocol
$ws.Range("D29").Value = "'37.46"
$ws.Range("E29").Value = "  -3.95%  "

# Row 30 - Filecoin
$ws.Range("E30").Value = "  +0.24%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'20.03"
$ws.Range("E31").Value = "  +0.18%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'150.65"
$ws.Range("E32").Value = "  -0.77%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +0.11%  "

# Row 34 - WEMIXToken
$ws.Range("D34").Value = "'2.64"
$ws.Range("E34").Value = "  +1.06%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -7.32%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -0.63%  "

# Row 37 - Stellar
$ws.Range("E37").Value = "  -0.64%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  +5.52%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "'15.06"
$ws.Range("E39").Value = "  +1.22%  "

# Row 40 - NEARProtocol
$ws.Range("E40").Value = "  -0.39%  "

# Row 41 - was RenderToken, now VeChain
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0304"
$ws.Range("E41").Value = "  +1.91%  "

# Row 42 - was VeChain, now RenderToken
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'3.80"
$ws.Range("E42").Value = "  -1.33%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.27%  "

# Row 44 - Maker
$ws.Range("D44").Value = "'1.839.85"
$ws.Range("E44").Value = "  +6.77%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  +14.46%  "

# Row 46 - BitcoinSV
$ws.Range("D46").Value = "'80.25"
$ws.Range("E46").Value = "  -4.15%  "

# Row 47 - Algorand
$ws.Range("E47").Value = "  +0.93%  "

# Row 48 - Aave
$ws.Range("D48").Value = "'99.43"
$ws.Range("E48").Value = "  -0.69%  "

# Row 49 - THORChain
$ws.Range("E49").Value = "  +0.38%  "

# Row 50 - ordi
$ws.Range("D50").Value = "'69.29"
$ws.Range("E50").Value = "  +0.16%  "

# Row 51 - MultiversX
$ws.Range("D51").Value = "'54.80"
$ws.Range("E51").Value = "  +0.87%  "
